# Apply the "fix: update backend interfaces" changes to the train.xlsx data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-03-09 06:05:12"

# Row 2
$ws.Range("D2").Value = $newTimestamp
$ws.Range("F2").Value = "http://49.234.6.241:5230/api/v1/memo/21"
$ws.Range("G2").Value = "/api/v1/memo/21"

# Row 3
$ws.Range("D3").Value = $newTimestamp
$ws.Range("M3").Value = 0.003

# Row 4
$ws.Range("D4").Value = $newTimestamp
$ws.Range("M4").Value = 0.003

# Row 5
$ws.Range("D5").Value = $newTimestamp
$ws.Range("M5").Value = 0.004

# Row 6
$ws.Range("D6").Value = $newTimestamp
$ws.Range("M6").Value = 0.003
$ws.Range("N6").Value = 0
$ws.Range("Q6").Value = $true

# Row 7
$ws.Range("B7").Value = 8
$ws.Range("D7").Value = $newTimestamp
$ws.Range("F7").Value = "http://49.234.6.241:5230/api/v1/memo/21"
$ws.Range("G7").Value = "/api/v1/memo/21"

# Row 8
$ws.Range("B8").Value = 9
$ws.Range("D8").Value = $newTimestamp

# Row 9
$ws.Range("B9").Value = 10
$ws.Range("D9").Value = $newTimestamp

# Row 10
$ws.Range("B10").Value = 11
$ws.Range("D10").Value = $newTimestamp
$ws.Range("F10").Value = "http://47.97.114.24:5230/api/v1/resource/16"
$ws.Range("G10").Value = "/api/v1/resource/16"
$ws.Range("N10").Value = 1
$ws.Range("Q10").Value = $false

# Row 11
$ws.Range("D11").Value = $newTimestamp
$ws.Range("F11").Value = "http://49.234.6.241:5230/api/v1/resource/16"
$ws.Range("G11").Value = "/api/v1/resource/16"
$ws.Range("M11").Value = 0.003
